$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 2762.5
$ws.Range("J29").Value = 5000
$ws.Range("L29").Value = 15000
$ws.Range("N29").Value = -15562

# Row 33
$ws.Range("H33").Value = 1015640.75
$ws.Range("I33").Value = 1232635.1
$ws.Range("K33").Value = 1232635.1
$ws.Range("M33").Value = -1232406.1

# Row 69
$ws.Range("H69").Value = 2222
$ws.Range("I69").Value = 2222
$ws.Range("K69").Value = 6666
$ws.Range("M69").Value = -5792

# Row 72
$ws.Range("H72").Value = 2222
$ws.Range("I72").Value = 2222
$ws.Range("K72").Value = 19998
$ws.Range("M72").Value = -15630

# Row 74
$ws.Range("H74").Value = 4414.143
$ws.Range("I74").Value = 3474.75
$ws.Range("J74").Value = 5666.6665
$ws.Range("K74").Value = 3474.75
$ws.Range("L74").Value = 5666.6665
$ws.Range("M74").Value = -2538.75
$ws.Range("N74").Value = -7538.6665

# Row 77
$ws.Range("H77").Value = 4414.143
$ws.Range("I77").Value = 3474.75
$ws.Range("J77").Value = 5666.6665
$ws.Range("K77").Value = 17373.75
$ws.Range("L77").Value = 28333.3325
$ws.Range("M77").Value = -12693.75
$ws.Range("N77").Value = -37693.3325

# Row 94
$ws.Range("H94").Value = 1030.5714
$ws.Range("I94").Value = 1030.5714
$ws.Range("K94").Value = 1030.5714
$ws.Range("M94").Value = -579.5714

# Row 132
$ws.Range("H132").Value = 4161.282
$ws.Range("I132").Value = 2820.6858
$ws.Range("K132").Value = 8462.057400000002
$ws.Range("M132").Value = -5932.057400000002

# Row 138
$ws.Range("H138").Value = 3629.1187
$ws.Range("I138").Value = 1574.8462
$ws.Range("J138").Value = 4209.674
$ws.Range("K138").Value = 4724.5386
$ws.Range("L138").Value = 12629.022
$ws.Range("M138").Value = 415.4614000000001
$ws.Range("N138").Value = -22909.022

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Range("H63").Value = 1656.5
$ws.Range("I63").Value = 1656.5
$ws.Range("K63").Value = 1656.5
$ws.Range("M63").Value = -970.5

# Row 66
$ws.Range("H66").Value = 1656.5
$ws.Range("I66").Value = 1656.5
$ws.Range("K66").Value = 8282.5
$ws.Range("M66").Value = -4850.5

# Row 74
$ws.Range("H74").Value = 2189.96
$ws.Range("I74").Value = 1934.1578
$ws.Range("K74").Value = 1934.1578
$ws.Range("M74").Value = -1060.1578

# Row 77
$ws.Range("H77").Value = 2189.96
$ws.Range("I77").Value = 1934.1578
$ws.Range("K77").Value = 9670.789000000001
$ws.Range("M77").Value = -5302.789000000001

# Row 110
$ws.Range("H110").Value = 2563.5
$ws.Range("I110").Value = 2206.1667
$ws.Range("J110").Value = 3099.5
$ws.Range("K110").Value = 2206.1667
$ws.Range("L110").Value = 3099.5
$ws.Range("M110").Value = -161.1667000000002
$ws.Range("N110").Value = -7189.5

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 1762.5
$ws.Range("I86").Value = 1543.2222
$ws.Range("K86").Value = 1543.2222
$ws.Range("M86").Value = -420.2221999999999

# Row 89
$ws.Range("H89").Value = 1762.5
$ws.Range("I89").Value = 1543.2222
$ws.Range("K89").Value = 7716.111
$ws.Range("M89").Value = -2100.111

# Row 105
$ws.Range("H105").Value = 5892748
$ws.Range("I105").Value = 7147779.5
$ws.Range("K105").Value = 7147779.5
$ws.Range("M105").Value = -7146032.5

# Row 134
$ws.Range("H134").Value = 3315.8262
$ws.Range("I134").Value = 2723.55
$ws.Range("J134").Value = 7264.3335
$ws.Range("K134").Value = 8170.650000000001
$ws.Range("L134").Value = 21793.0005
$ws.Range("M134").Value = -5635.650000000001
$ws.Range("N134").Value = -26863.0005

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 3007.9092
$ws.Range("J16").Value = 2985
$ws.Range("L16").Value = 2985
$ws.Range("N16").Value = -3559

# Row 31
$ws.Range("H31").Value = 3709.25
$ws.Range("I31").Value = 2837.2727
$ws.Range("K31").Value = 2837.2727
$ws.Range("M31").Value = -2542.2727

# Row 34
$ws.Range("H34").Value = 3709.25
$ws.Range("I34").Value = 2837.2727
$ws.Range("K34").Value = 2837.2727
$ws.Range("M34").Value = -2635.2727

# Row 99
$ws.Range("H99").Value = 71888.03999999999
$ws.Range("J99").Value = 1997.25
$ws.Range("L99").Value = 1997.25
$ws.Range("N99").Value = -4993.25

# Row 105
$ws.Range("H105").Value = 1889.6
$ws.Range("I105").Value = 1862.8334
$ws.Range("J105").Value = 1996.6666
$ws.Range("K105").Value = 1862.8334
$ws.Range("L105").Value = 1996.6666
$ws.Range("M105").Value = -115.8334
$ws.Range("N105").Value = -5490.6666

# Row 113
$ws.Range("H113").Value = 3007.9092
$ws.Range("J113").Value = 2985
$ws.Range("L113").Value = 2985
$ws.Range("N113").Value = -7325

# Row 126
$ws.Range("H126").Value = 71888.03999999999
$ws.Range("J126").Value = 1997.25
$ws.Range("L126").Value = 5991.75
$ws.Range("N126").Value = -10931.75

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 1002445
$ws.Range("I68").Value = 2356.7083
$ws.Range("K68").Value = 7070.124899999999
$ws.Range("M68").Value = -6259.124899999999

# Row 71
$ws.Range("H71").Value = 1002445
$ws.Range("I71").Value = 2356.7083
$ws.Range("K71").Value = 21210.3747
$ws.Range("M71").Value = -17154.3747

# Row 109
$ws.Range("H109").Value = 3215.8235
$ws.Range("I109").Value = 808.75
$ws.Range("K109").Value = 2426.25
$ws.Range("M109").Value = -1386.25

# Row 137
$ws.Range("H137").Value = 5559176.5
$ws.Range("J137").Value = 5847.3335
$ws.Range("L137").Value = 17542.0005
$ws.Range("N137").Value = -27742.0005

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 66676772
$ws.Range("I80").Value = 90920824
$ws.Range("K80").Value = 90920824
$ws.Range("M80").Value = -90919826

# Row 83
$ws.Range("H83").Value = 66676772
$ws.Range("I83").Value = 90920824
$ws.Range("K83").Value = 454604120
$ws.Range("M83").Value = -454599128

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 1845.6923
$ws.Range("I46").Value = 1408.2
$ws.Range("K46").Value = 1408.2
$ws.Range("M46").Value = -1220.2

# Row 55
$ws.Range("H55").Value = 280.1
$ws.Range("I55").Value = 228
$ws.Range("J55").Value = 384.3
$ws.Range("K55").Value = 228
$ws.Range("L55").Value = 384.3
$ws.Range("M55").Value = -55
$ws.Range("N55").Value = -730.3

# Row 82
$ws.Range("H82").Value = 1132.2
$ws.Range("I82").Value = 969.4
$ws.Range("J82").Value = 1457.8
$ws.Range("K82").Value = 969.4
$ws.Range("L82").Value = 1457.8
$ws.Range("M82").Value = -608.4
$ws.Range("N82").Value = -2179.8

# Row 85
$ws.Range("H85").Value = 1132.2
$ws.Range("I85").Value = 969.4
$ws.Range("J85").Value = 1457.8
$ws.Range("K85").Value = 969.4
$ws.Range("L85").Value = 1457.8
$ws.Range("M85").Value = 278.6
$ws.Range("N85").Value = -3953.8
